$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: measure/dimension metadata - municipio & aragon dimensions collapsed
# into a single shared "sdmx-dimension:refArea" dimension.
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# Row 3: municipio column role changes from "medida" (measure) to "dim" (dimension)
$ws.Range("C3").Value = "dim"

# Row 4: URI type columns - municipio becomes a URI column, aragon becomes "URI-Comunidad"
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("E4").Value = "URI-Comunidad"

# Row 5: drop the now-unused mapping-aragon.xlsx reference (was in E5)
$ws.Range("E5").Clear()
